$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new quarter columns before column D (shifts D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number formats/styles from column F (the old column D data) into new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new quarter values (6/30/2018 and 9/30/2018 columns)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 26000
$ws.Range("E8").Value = 30100
$ws.Range("D9").Value = 6800
$ws.Range("E9").Value = 6400
$ws.Range("D10").Value = 19200
$ws.Range("E10").Value = 23700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 500
$ws.Range("D15").Value = 900
$ws.Range("E15").Value = 900
$ws.Range("D17").Value = 12200
$ws.Range("E17").Value = 13000
$ws.Range("D18").Value = 13800
$ws.Range("E18").Value = 17100
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 14700
$ws.Range("E21").Value = 17900
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 13800
$ws.Range("E23").Value = 17100
$ws.Range("D24").Value = -100
$ws.Range("E24").Value = -800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 13900
$ws.Range("E26").Value = 17900
$ws.Range("D27").Value = 13100
$ws.Range("E27").Value = 14900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 13100
$ws.Range("E33").Value = 14900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 13100
$ws.Range("E35").Value = 14900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 32000
$ws.Range("E41").Value = 25000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 7000
$ws.Range("E43").Value = 7900
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 873300
$ws.Range("E47").Value = 895300
$ws.Range("D48").Value = 64600
$ws.Range("E48").Value = 65300
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1700
$ws.Range("E52").Value = 1200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 982700
$ws.Range("E54").Value = 1001200
$ws.Range("D57").Value = 7500
$ws.Range("E57").Value = 7800
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 7600
$ws.Range("E59").Value = 9700
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 568800
$ws.Range("E61").Value = 600900
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 583900
$ws.Range("E66").Value = 618400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 94700
$ws.Range("E70").Value = 94500
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("E72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 304100
$ws.Range("E76").Value = 288300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 13100
$ws.Range("E81").Value = 14900
$ws.Range("D83").Value = 900
$ws.Range("E83").Value = 900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 7700
$ws.Range("E89").Value = 10100
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 41700
$ws.Range("E94").Value = 22400
$ws.Range("D96").Value = -10400
$ws.Range("E96").Value = -8300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -41700
$ws.Range("E100").Value = -34700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 7600
$ws.Range("E102").Value = -2100
